$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 379, shifting existing rows 379:467 down to 380:468
$ws.Rows.Item(379).EntireRow.Insert()

# Populate the newly inserted row 379 with the new data record
$ws.Cells.Item(379, 1).Value = 4
$ws.Cells.Item(379, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(379, 3).Value = "Los Lagos"
$ws.Cells.Item(379, 4).Value = 45244
$ws.Cells.Item(379, 5).Value = 10
$ws.Cells.Item(379, 6).Value = 100112044
$ws.Cells.Item(379, 7).Value = "Perejil"
$ws.Cells.Item(379, 8).Value = "Sin especificar"
$ws.Cells.Item(379, 9).Value = "Primera"
$ws.Cells.Item(379, 10).Value = 160
$ws.Cells.Item(379, 11).Value = 8000
$ws.Cells.Item(379, 12).Value = 8000
$ws.Cells.Item(379, 13).Value = 8000
$ws.Cells.Item(379, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(379, 15).Value = "Región Metropolitana"
$ws.Cells.Item(379, 16).Value = 2667
$ws.Cells.Item(379, 17).Value = 3
$ws.Cells.Item(379, 18).Value = "Hortaliza"

# Match the date style used by column D elsewhere (D379 should carry style index 2, same as D380 etc.)
$ws.Cells.Item(379, 4).NumberFormat = $ws.Cells.Item(380, 4).NumberFormat
